$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to Text format so numeric-looking strings (e.g. "1.002") are not
# auto-converted to numbers, preserving the original inlineStr/text semantics.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.846.09"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").Value = "1.893.03"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "0.7562"
$ws.Range("E5").Value = "  +3.06%  "

$ws.Range("D6").Value = "239.72"
$ws.Range("E6").Value = "  -1.69%  "

$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").Value = "0.3032"
$ws.Range("E8").Value = "  -3.16%  "

$ws.Range("D9").Value = "25.25"
$ws.Range("E9").Value = "  -6.03%  "

$ws.Range("D10").Value = "0.06798"
$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("D11").Value = "0.07962"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.891.17"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7394"
$ws.Range("E13").Value = "  -4.94%  "

$ws.Range("D14").Value = "5.150"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").Value = "90.73"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "29.878.28"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").Value = "5.924"
$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("D19").Value = "241.52"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "0.000007680"
$ws.Range("E20").Value = "  -1.45%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").Value = "6.884"
$ws.Range("E23").Value = "  +2.00%  "

$ws.Range("D24").Value = "9.192"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").Value = "165.21"
$ws.Range("E25").Value = "  -0.46%  "

$ws.Range("D26").Value = "18.60"
$ws.Range("E26").Value = "  -2.35%  "

$ws.Range("D27").Value = "0.1268"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").Value = "2.005"
$ws.Range("E28").Value = "  -3.90%  "

$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  +2.76%  "

$ws.Range("D30").Value = "1.513"
$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("D31").Value = "4.232"
$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("D32").Value = "4.001"
$ws.Range("E32").Value = "  -1.99%  "

$ws.Range("D33").Value = "0.05275"
$ws.Range("E33").Value = "  +2.22%  "

$ws.Range("D34").Value = "1.236"
$ws.Range("E34").Value = "  -3.74%  "

$ws.Range("D35").Value = "0.7176"
$ws.Range("E35").Value = "  -3.44%  "

$ws.Range("D36").Value = "2.716"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("D37").Value = "0.01902"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").Value = "2.788"
$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").Value = "6.130"
$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").Value = "0.4370"
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("D41").Value = "71.53"
$ws.Range("E41").Value = "  -4.21%  "

$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "1.871"
$ws.Range("E43").Value = "  -3.20%  "

$ws.Range("D44").Value = "0.8248"

$ws.Range("D45").Value = "100.53"
$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("D46").Value = "9.785"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "7.469"
$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("D48").Value = "2.049.88"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").Value = "36.09"
$ws.Range("E49").Value = "  -3.85%  "

$ws.Range("D50").Value = "0.05959"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("D51").Value = "1.463"
$ws.Range("E51").Value = "  +0.47%  "

# Reset style back to Normal so no stray cell style (s attribute) is introduced
# by the temporary Text number-format applied above.
$ws.Range("D2:E51").Style = "Normal"
